# Remove the "Appendix: Quick prototype" section: the Heading2 title,
# the three "Figure: PDF page NN" captions, and their embedded
# screenshot images -- everything that sits between the earlier
# "Appendix: LinksHYPERLINK ..." paragraph and the final Heading2
# "Appendix: Links" / "ADO Epic Link" paragraphs.

$d = $word.ActiveDocument

$startHeading = "Appendix: Quick prototype"
$endHeading   = "Appendix: Links"

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)

    if ($null -eq $startPara) {
        if ($text -eq $startHeading) {
            $startPara = $p
        }
        continue
    }

    if ($text -eq $endHeading) {
        $endPara = $p
        break
    }
}

if ($null -ne $startPara -and $null -ne $endPara) {
    # Delete from the start of the "Quick prototype" heading up to (but
    # not including) the start of the final "Appendix: Links" heading,
    # taking the appendix heading, figure captions, and inline images
    # with it while leaving the surrounding paragraphs intact.
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.Start)
    $rng.Delete()
}
